# Update CN KeyItem.xlsx:
#  - Insert two new key items ("552" Bubbling Mysterious Vial - corruption scent,
#    "553" Bubbling Mysterious Vial - detergent scent) right after "550"
#    (Bracelet of Forgiveness), which pushes the existing 590/600/601/610 rows
#    down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as text even
# when it looks like a number (e.g. "552"), without leaving a NumberFormat
# override behind on the target cell itself.
function Set-TextValue($range, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $scratch.Clear() | Out-Null
}

# Helper: write an entire row (hashtable keyed by column letter) into the sheet.
function Set-RowData($rowNum, $data) {
    foreach ($col in @("A","B","C","D","E","F","G","H")) {
        $cellRange = $ws.Range("$col$rowNum")
        $value = $data[$col]
        if ($col -eq "A") {
            Set-TextValue $cellRange $value
        } else {
            $cellRange.Value = $value
        }
    }
}

# Insert two fresh rows at 11-12 (shifts old rows 11-14 down to 13-16),
# copying the formatting of the row right above so the new rows look the same
# as their neighbours.
$ws.Range("A10:H10").Copy() | Out-Null
$ws.Rows("11:12").Insert(-4121) | Out-Null   # xlShiftDown

# New row 11: id 552 - corruption-scented vial
$row11 = @{
    A = "552"
    B = "EA 23.207"
    C = "装有起泡液体的神秘瓶子"
    D = "Bubbling Mysterious Vial"
    E = "泡立つ謎の瓶"
    F = "这是装有起泡的透明液体的瓶子。它散发出了污秽的臭味。"
    G = "A vial filled with a bubbling, clear liquid. It reeks of corruption."
    H = "泡立つ透明な液体が入った瓶だ。とても穢れた匂いがする。"
}
Set-RowData 11 $row11

# New row 12: id 553 - detergent-scented vial
$row12 = @{
    A = "553"
    B = "EA 23.207"
    C = "装有起泡液体的神秘瓶子"
    D = "Bubbling Mysterious Vial"
    E = "泡立つ謎の瓶"
    F = "这是装有起泡的透明液体的瓶子。它散发出了清洗剂的味道。"
    G = "A vial filled with a bubbling, clear liquid. It smells like detergent."
    H = "泡立つ透明な液体が入った瓶だ。洗剤の匂いがする。"
}
Set-RowData 12 $row12

Write-Host "Done. A11=$($ws.Range('A11').Value()) A12=$($ws.Range('A12').Value()) A13=$($ws.Range('A13').Value()) A16=$($ws.Range('A16').Value())"
